$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2022471910112359
$ws.Range("C2").Value = 0.5252808988764045
$ws.Range("J2").Value = 0.03089887640449438
$ws.Range("P2").Value = 0.1348314606741573
$ws.Range("S2").Value = 0.1067415730337079
$ws.Range("B3").Value = 0.0154639175257732
$ws.Range("C3").Value = 0.02577319587628866
$ws.Range("J3").Value = 0.005154639175257732
$ws.Range("P3").Value = 0.7577319587628866
$ws.Range("S3").Value = 0.1958762886597938
$ws.Range("J4").Value = 0.01639344262295082
$ws.Range("P4").Value = 0.7049180327868853
$ws.Range("S4").Value = 0.2786885245901639
$ws.Range("B6").Value = 0.06143344709897611
$ws.Range("D6").Value = 0.0136518771331058
$ws.Range("F6").Value = 0.07849829351535836
$ws.Range("J6").Value = 0.2389078498293515
$ws.Range("O6").Value = 0.0273037542662116
$ws.Range("Q6").Value = 0.1467576791808874
$ws.Range("R6").Value = 0.06825938566552901
$ws.Range("S6").Value = 0.3651877133105802
$ws.Range("B7").Value = 0.09411764705882353
$ws.Range("D7").Value = 0.04117647058823529
$ws.Range("F7").Value = 0.07058823529411765
$ws.Range("J7").Value = 0.1411764705882353
$ws.Range("O7").Value = 0.02941176470588235
$ws.Range("Q7").Value = 0.1411764705882353
$ws.Range("R7").Value = 0.06470588235294118
$ws.Range("S7").Value = 0.4176470588235294
$ws.Range("B8").Value = 0.112701252236136
$ws.Range("D8").Value = 0.02862254025044723
$ws.Range("E8").Value = 0.003577817531305903
$ws.Range("F8").Value = 0.09302325581395349
$ws.Range("J8").Value = 0.09481216457960644
$ws.Range("O8").Value = 0.02504472271914132
$ws.Range("Q8").Value = 0.1359570661896243
$ws.Range("R8").Value = 0.08228980322003578
$ws.Range("S8").Value = 0.4239713774597496
$ws.Range("B9").Value = 0.06854838709677419
$ws.Range("D9").Value = 0.03225806451612903
$ws.Range("F9").Value = 0.08064516129032258
$ws.Range("J9").Value = 0.08870967741935484
$ws.Range("O9").Value = 0.01612903225806452
$ws.Range("Q9").Value = 0.1491935483870968
$ws.Range("R9").Value = 0.09274193548387097
$ws.Range("S9").Value = 0.4717741935483871
$ws.Range("B10").Value = 0.1357615894039735
$ws.Range("D10").Value = 0.02235099337748344
$ws.Range("E10").Value = 0.002483443708609272
$ws.Range("F10").Value = 0.07864238410596026
$ws.Range("J10").Value = 0.07533112582781457
$ws.Range("O10").Value = 0.02649006622516556
$ws.Range("Q10").Value = 0.2044701986754967
$ws.Range("R10").Value = 0.06705298013245033
$ws.Range("S10").Value = 0.3874172185430463
$ws.Range("F11").Value = 0.004219409282700422
$ws.Range("G11").Value = 0.1476793248945148
$ws.Range("J11").Value = 0.05907172995780591
$ws.Range("K11").Value = 0.1687763713080169
$ws.Range("L11").Value = 0.5991561181434599
$ws.Range("S11").Value = 0.02109704641350211
$ws.Range("G12").Value = 0.7432432432432432
$ws.Range("J12").Value = 0.1689189189189189
$ws.Range("K12").Value = 0.01351351351351351
$ws.Range("L12").Value = 0.02702702702702703
$ws.Range("S12").Value = 0.0472972972972973
$ws.Range("F13").Value = 0.02
$ws.Range("G13").Value = 0.7
$ws.Range("J13").Value = 0.26
$ws.Range("S13").Value = 0.02
$ws.Range("F15").Value = 0.03249097472924187
$ws.Range("H15").Value = 0.1913357400722022
$ws.Range("I15").Value = 0.07581227436823104
$ws.Range("J15").Value = 0.3465703971119133
$ws.Range("K15").Value = 0.02527075812274368
$ws.Range("M15").Value = 0.02166064981949458
$ws.Range("O15").Value = 0.06137184115523465
$ws.Range("S15").Value = 0.2454873646209386
$ws.Range("F16").Value = 0.04680851063829787
$ws.Range("H16").Value = 0.2042553191489362
$ws.Range("I16").Value = 0.0851063829787234
$ws.Range("J16").Value = 0.3787234042553191
$ws.Range("K16").Value = 0.09361702127659574
$ws.Range("M16").Value = 0.008510638297872341
$ws.Range("O16").Value = 0.06382978723404255
$ws.Range("S16").Value = 0.1191489361702128
$ws.Range("F17").Value = 0.02132701421800948
$ws.Range("H17").Value = 0.2061611374407583
$ws.Range("I17").Value = 0.0995260663507109
$ws.Range("J17").Value = 0.4028436018957346
$ws.Range("K17").Value = 0.08056872037914692
$ws.Range("M17").Value = 0.02132701421800948
$ws.Range("O17").Value = 0.05687203791469194
$ws.Range("S17").Value = 0.1113744075829384
$ws.Range("F18").Value = 0.03846153846153846
$ws.Range("H18").Value = 0.2197802197802198
$ws.Range("I18").Value = 0.1318681318681319
$ws.Range("J18").Value = 0.3296703296703297
$ws.Range("K18").Value = 0.0989010989010989
$ws.Range("M18").Value = 0.01648351648351648
$ws.Range("N18").Value = 0.005494505494505495
$ws.Range("O18").Value = 0.08241758241758242
$ws.Range("S18").Value = 0.07692307692307693
$ws.Range("F19").Value = 0.01918976545842218
$ws.Range("H19").Value = 0.2380952380952381
$ws.Range("I19").Value = 0.09950248756218906
$ws.Range("J19").Value = 0.3461265103056148
$ws.Range("K19").Value = 0.07818052594171997
$ws.Range("M19").Value = 0.02203269367448472
$ws.Range("N19").Value = 0.001421464108031272
$ws.Range("O19").Value = 0.08102345415778252
$ws.Range("S19").Value = 0.1144278606965174
